$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the two newly-filled participants (206 = script "06", 207 = script "07").
# Columns: A=Group, B=Script, C=Event, D=Dispersion_School, E=Confident_Factor, F=Disruption_Factor
$rows = @(
    @(47, 206, "06", "clicking pen",     4, 10, 10),
    @(48, 206, "06", "locking at phone", 2, 10, 8),
    @(49, 206, "06", "drumming",         2, 10, 10),
    @(50, 206, "06", "drawing",          2, 10, 10),
    @(51, 206, "06", "head on table",    4, 10, 7),
    @(52, 206, "06", "whispering",       5, 9,  6),
    @(53, 206, "06", "snipping",         3, 9,  4),
    @(54, 206, "06", "heckling",         2, 10, 6),
    @(55, 206, "06", "chatting",         1, 10, 7),
    @(56, 207, "07", "locking at phone", 4, 8,  7),
    @(57, 207, "07", "drawing",          1, 10, 8),
    @(58, 207, "07", "clicking pen",     10, 10, 2),
    @(59, 207, "07", "whispering",       5, 5,  7),
    @(60, 207, "07", "drumming",         5, 8,  8),
    @(61, 207, "07", "heckling",         7, 6,  3),
    @(62, 207, "07", "head on table",    6, 5,  7),
    @(63, 207, "07", "snipping",         0, 9,  2),
    @(64, 207, "07", "chatting",         7, 7,  4)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}

# Update the sheet's active-cell selection to match the post-edit state.
$ws.Range("D67").Select()
